# Add a new "Grid" process set row to the VEDA_Sets-Proc sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEDA_Sets-Proc")

# Populate in the same order the shared-string table was built by the
# original author (SetName first, then PSET_SET, then PSET_PN).
$ws.Range("F21").Value = "Grid"
$ws.Range("A21").Value = "IRE"
$ws.Range("B21").Value = "g[_]*"

$ws.Range("B21").Select() | Out-Null
